# added 4wk low sales check
# Update the "Forecast Comparison" sheet with refreshed forecast figures
# (MyForecast, Inventory Coverage, Stockout Risk, Reorder Urgency,
# Seasonality Index) and roll the new totals into the "Summary" sheet.

$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# Row layout: D=MyForecast, H=Inventory Coverage, I=Stockout Risk,
#             J=Reorder Urgency, L=Seasonality Index
$rows = @(
    @{ Row = 2;  D = 6; H = 15.33; L = 1 },
    @{ Row = 3;  D = 5; H = 17.2;  L = 0.8100000000000001 },
    @{ Row = 4;  D = 6; H = 13.5;  L = 0.84 },
    @{ Row = 5;  D = 7; H = 10.71; L = 1.19 },
    @{ Row = 6;  D = 7; H = 9.710000000000001;  L = 1.03 },
    @{ Row = 7;  D = 7; H = 8.710000000000001;  L = 0.9 },
    @{ Row = 8;  H = 7.71; L = 1.04 },
    @{ Row = 9;  H = 6.71; L = 1.05 },
    @{ Row = 10; D = 7; H = 5.71; L = 1.15 },
    @{ Row = 11; D = 8; H = 4.12; L = 0.93 },
    @{ Row = 12; D = 8; H = 3.12; L = 0.93 },
    @{ Row = 13; H = 2.43; L = 0.88 },
    @{ Row = 14; D = 8; H = 1.25; L = 1.11 },
    @{ Row = 15; D = 8; H = 0.25; I = "High"; J = "Urgent"; L = 1.11 },
    @{ Row = 16; D = 8; H = 0;    I = "High"; J = "Urgent"; L = 0.8100000000000001 },
    @{ Row = 17; D = 8; H = 0;    I = "High"; J = "Urgent"; L = 0.89 }
)

foreach ($r in $rows) {
    $row = $r.Row
    if ($r.ContainsKey("D")) { $wsForecast.Range("D$row").Value = $r.D }
    if ($r.ContainsKey("H")) { $wsForecast.Range("H$row").Value = $r.H }
    if ($r.ContainsKey("I")) { $wsForecast.Range("I$row").Value = $r.I }
    if ($r.ContainsKey("J")) { $wsForecast.Range("J$row").Value = $r.J }
    if ($r.ContainsKey("L")) { $wsForecast.Range("L$row").Value = $r.L }
}

# Refresh Summary sheet totals derived from the updated MyForecast column.
# These cells store their numbers as TEXT (not numeric values), so force a
# text number-format before writing, otherwise Excel auto-coerces a
# numeric-looking string into a real number. Restore the default "Normal"
# style afterwards so no stray per-cell style survives the round trip.
$summaryTextCells = @(
    @{ Cell = "B9";  Value = "114" },
    @{ Cell = "B10"; Value = "52" },
    @{ Cell = "B11"; Value = "24" },
    @{ Cell = "B14"; Value = "5" }
)

foreach ($c in $summaryTextCells) {
    $rng = $wsSummary.Range($c.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $c.Value
    $rng.Style = "Normal"
}
